# Updates generated data values (want-to-go counts / lowest prices / sale
# status) across the four sheets, matching the "output generated" refresh
# commit. Values are written with .Value so Excel stores them as the
# appropriate cell type (number or string).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3550
$ws1.Range("F5").Value  = 8310
$ws1.Range("F7").Value  = 119
$ws1.Range("F10").Value = 93
$ws1.Range("F11").Value = 70
$ws1.Range("F14").Value = 7325
$ws1.Range("G15").Value = "不可售"
$ws1.Range("F16").Value = 7611
$ws1.Range("F18").Value = 57472
$ws1.Range("F19").Value = 57472
$ws1.Range("F20").Value = 4723
$ws1.Range("F28").Value = 5284
$ws1.Range("F30").Value = 104
$ws1.Range("F33").Value = 1341
$ws1.Range("F34").Value = 1790
$ws1.Range("F37").Value = 225
$ws1.Range("F41").Value = 40
$ws1.Range("F43").Value = 251
$ws1.Range("F44").Value = 126
$ws1.Range("F47").Value = 193
$ws1.Range("F49").Value = 56

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 195
$ws2.Range("F6").Value  = 134
$ws2.Range("G6").Value  = 280
$ws2.Range("F10").Value = 7591
$ws2.Range("F15").Value = 2
$ws2.Range("F24").Value = 30
$ws2.Range("F26").Value = 1
$ws2.Range("F32").Value = 80
$ws2.Range("F37").Value = 50
$ws2.Range("F48").Value = 275

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value  = 2356
$ws3.Range("F9").Value  = 9424
$ws3.Range("F10").Value = 1747
$ws3.Range("F15").Value = 263
$ws3.Range("F16").Value = 2304
$ws3.Range("F18").Value = 488

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2356
$ws4.Range("F4").Value  = 8310
$ws4.Range("F6").Value  = 1747
$ws4.Range("F8").Value  = 119
$ws4.Range("F9").Value  = 2304
$ws4.Range("F10").Value = 70
$ws4.Range("F11").Value = 7611
$ws4.Range("F12").Value = 57472
$ws4.Range("F13").Value = 195
$ws4.Range("F16").Value = 4723
$ws4.Range("F23").Value = 5284
$ws4.Range("F25").Value = 104
$ws4.Range("F26").Value = 44
$ws4.Range("F28").Value = 1341
$ws4.Range("F29").Value = 1791
$ws4.Range("F31").Value = 488
$ws4.Range("F36").Value = 225
$ws4.Range("F41").Value = 251
$ws4.Range("F43").Value = 30
$ws4.Range("F46").Value = 56
$ws4.Range("F51").Value = 275
